$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 14: change the "category tag" (column H) from "normal" to "special"
$ws.Range("H14").Value = "special"

# Add a new value in column K (是否常驻 / always present) for row 14
$ws.Range("K14").Value = "yes"

# Update the active selection to K15, matching the saved view state
$ws.Activate()
$ws.Range("K15").Select()
